$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 0.99999999440790988
$ws.Cells.Item(2, 1).Value = 0.99546598295946542
$ws.Cells.Item(3, 1).Value = 0.97666940519947609
$ws.Cells.Item(4, 1).Value = 0.96865985497642515
$ws.Cells.Item(5, 1).Value = 0.96111935027633644
$ws.Cells.Item(6, 1).Value = 0.9451814911470352
$ws.Cells.Item(7, 1).Value = 0.9424972618140921
$ws.Cells.Item(8, 1).Value = 0.93614667863048118
$ws.Cells.Item(9, 1).Value = 0.93066294666093596
$ws.Cells.Item(10, 1).Value = 0.92609014529050981
$ws.Cells.Item(11, 1).Value = 0.92524863269168667
$ws.Cells.Item(12, 1).Value = 0.92390429928210738
$ws.Cells.Item(13, 1).Value = 0.92200757176386061
$ws.Cells.Item(14, 1).Value = 0.91784034897956446
$ws.Cells.Item(15, 1).Value = 0.91524892235444433
$ws.Cells.Item(16, 1).Value = 0.91274247253012319
$ws.Cells.Item(17, 1).Value = 0.9090345861541651
$ws.Cells.Item(18, 1).Value = 0.90792567320173578
$ws.Cells.Item(19, 1).Value = 0.99385092014638832
$ws.Cells.Item(20, 1).Value = 0.98158949469708301
$ws.Cells.Item(21, 1).Value = 0.97886005809070953
$ws.Cells.Item(22, 1).Value = 0.97536096843369802
$ws.Cells.Item(23, 1).Value = 0.96963404023025346
$ws.Cells.Item(24, 1).Value = 0.95661258634848789
$ws.Cells.Item(25, 1).Value = 0.95015549373373087
$ws.Cells.Item(26, 1).Value = 0.94402465844226369
$ws.Cells.Item(27, 1).Value = 0.94266400955218754
$ws.Cells.Item(28, 1).Value = 0.93842907913794504
$ws.Cells.Item(29, 1).Value = 0.93600103045614702
$ws.Cells.Item(30, 1).Value = 0.93547151676775187
$ws.Cells.Item(31, 1).Value = 0.93336425191089389
$ws.Cells.Item(32, 1).Value = 0.93327823783730057
$ws.Cells.Item(33, 1).Value = 0.93275823121931456
